$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 265
$ws1.Range("F4").Value = 2673
$ws1.Range("F6").Value = 571

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 265
$ws4.Range("F6").Value = 2673
$ws4.Range("F8").Value = 571
